$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C are plain text (coin name / URL) -- safe to set directly.
# Columns D and E contain numeric-looking / formatted strings (e.g. "1.00", "0.999",
# "  +4.45%  ") that Excel would otherwise coerce into numbers and mangle
# (dropping trailing zeros, trimming padding spaces, etc). Force each such cell to
# Text format before assignment, then reset the style back to Normal so we don't
# leave a stray NumberFormat on the cell.

# D2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.542.26"
$ws.Range("D2").Style = "Normal"

# E2
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.89%  "
$ws.Range("E2").Style = "Normal"

# D3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.434.59"
$ws.Range("D3").Style = "Normal"

# E3
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +6.24%  "
$ws.Range("E3").Style = "Normal"

# D4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"

# D5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.19"
$ws.Range("D5").Style = "Normal"

# E5
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.22%  "
$ws.Range("E5").Style = "Normal"

# D6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.31"
$ws.Range("D6").Style = "Normal"

# E6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +8.19%  "
$ws.Range("E6").Style = "Normal"

# E7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E7").Style = "Normal"

# D8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.583"
$ws.Range("D8").Style = "Normal"

# E8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("E8").Style = "Normal"

# D9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.433.10"
$ws.Range("D9").Style = "Normal"

# E9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.28%  "
$ws.Range("E9").Style = "Normal"

# E10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.57%  "
$ws.Range("E10").Style = "Normal"

# D11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.76"
$ws.Range("D11").Style = "Normal"

# E11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.95%  "
$ws.Range("E11").Style = "Normal"

# E12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("E12").Style = "Normal"

# E13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.77%  "
$ws.Range("E13").Style = "Normal"

# D14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.32"
$ws.Range("D14").Style = "Normal"

# E14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +13.76%  "
$ws.Range("E14").Style = "Normal"

# D15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.865.54"
$ws.Range("D15").Style = "Normal"

# E15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +6.26%  "
$ws.Range("E15").Style = "Normal"

# D16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.404.42"
$ws.Range("D16").Style = "Normal"

# E16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.70%  "
$ws.Range("E16").Style = "Normal"

# E17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +8.14%  "
$ws.Range("E17").Style = "Normal"

# D18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.430.59"
$ws.Range("D18").Style = "Normal"

# E18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.43%  "
$ws.Range("E18").Style = "Normal"

# E19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.09%  "
$ws.Range("E19").Style = "Normal"

# D20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.62"
$ws.Range("D20").Style = "Normal"

# E20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +12.06%  "
$ws.Range("E20").Style = "Normal"

# E21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.94%  "
$ws.Range("E21").Style = "Normal"

# D22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.82"
$ws.Range("D22").Style = "Normal"

# E22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.44%  "
$ws.Range("E22").Style = "Normal"

# D23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"

# E23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E23").Style = "Normal"

# D24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.32"
$ws.Range("D24").Style = "Normal"

# E24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.85%  "
$ws.Range("E24").Style = "Normal"

# E25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("E25").Style = "Normal"

# D26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"

# E26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E26").Style = "Normal"

# E27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +15.38%  "
$ws.Range("E27").Style = "Normal"

# D28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.20"
$ws.Range("D28").Style = "Normal"

# E28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.36%  "
$ws.Range("E28").Style = "Normal"

# E29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +15.22%  "
$ws.Range("E29").Style = "Normal"

# D30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0790"
$ws.Range("D30").Style = "Normal"

# E30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +9.70%  "
$ws.Range("E30").Style = "Normal"

# E31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.42%  "
$ws.Range("E31").Style = "Normal"

# E32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +11.04%  "
$ws.Range("E32").Style = "Normal"

# D33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "172.00"
$ws.Range("D33").Style = "Normal"

# E33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("E33").Style = "Normal"

# D34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.45"
$ws.Range("D34").Style = "Normal"

# E34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.69%  "
$ws.Range("E34").Style = "Normal"

# D35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.397"
$ws.Range("D35").Style = "Normal"

# E35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.79%  "
$ws.Range("E35").Style = "Normal"

# B36
$ws.Range("B36").Value = "Bittensor"

# C36
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"

# D36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "376.96"
$ws.Range("D36").Style = "Normal"

# E36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +19.81%  "
$ws.Range("E36").Style = "Normal"

# B37
$ws.Range("B37").Value = "EthereumClassic"

# C37
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"

# D37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.60"
$ws.Range("D37").Style = "Normal"

# E37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.85%  "
$ws.Range("E37").Style = "Normal"

# B38
$ws.Range("B38").Value = "NEARProtocol"

# C38
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"

# D38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.51"
$ws.Range("D38").Style = "Normal"

# E38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +13.50%  "
$ws.Range("E38").Style = "Normal"

# E39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E39").Style = "Normal"

# D40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").Style = "Normal"

# E40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E40").Style = "Normal"

# E41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +12.60%  "
$ws.Range("E41").Style = "Normal"

# D42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.15"
$ws.Range("D42").Style = "Normal"

# E42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.61%  "
$ws.Range("E42").Style = "Normal"

# D43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "145.90"
$ws.Range("D43").Style = "Normal"

# E43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +7.50%  "
$ws.Range("E43").Style = "Normal"

# E44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +8.04%  "
$ws.Range("E44").Style = "Normal"

# D45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.88"
$ws.Range("D45").Style = "Normal"

# E45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +12.37%  "
$ws.Range("E45").Style = "Normal"

# D46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.591"
$ws.Range("D46").Style = "Normal"

# E46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.81%  "
$ws.Range("E46").Style = "Normal"

# E47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("E47").Style = "Normal"

# D48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0520"
$ws.Range("D48").Style = "Normal"

# E48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.86%  "
$ws.Range("E48").Style = "Normal"

# D49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0223"
$ws.Range("D49").Style = "Normal"

# E49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.03%  "
$ws.Range("E49").Style = "Normal"

# D50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.93"
$ws.Range("D50").Style = "Normal"

# E50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +7.74%  "
$ws.Range("E50").Style = "Normal"

# D51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0219"
$ws.Range("D51").Style = "Normal"

# E51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.48%  "
$ws.Range("E51").Style = "Normal"

